# Update natmi LR-pair data rows (ligand-expressing / receptor-expressing cell counts and
# all derived expression / specificity / weight statistics) per revised Natmi run
# (following Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2 = @{ "E" = 3; "G" = 6.201571333333333; "H" = 18.604714; "I" = 0.05221490529364391; "J" = 0.07406232529850043; "K" = 3; "M" = 1.253707; "N" = 3.761121; "O" = 0.07760709912293902; "P" = 0.1037097015912075; "Q" = 7.774953391599333; "R" = 69.97458052439399; "S" = 0.004052247330818696; "T" = 0.007680981655858415 }
    3 = @{ "E" = 3; "G" = 6.201571333333333; "H" = 18.604714; "I" = 0.05221490529364391; "J" = 0.07406232529850043; "K" = 3; "M" = 1.168375666666667; "N" = 3.505127; "O" = 0.07232491018701337; "P" = 0.09665088552303537; "Q" = 7.245765040964222; "R" = 65.21188536867798; "S" = 0.003776438335786204; "T" = 0.007158189323995172 }
    4 = @{ "E" = 3; "G" = 6.201571333333333; "H" = 18.604714; "I" = 0.05221490529364391; "J" = 0.07406232529850043; "K" = 3; "M" = 0.6422343333333334; "N" = 1.926703; "O" = 0.03975565548182683; "P" = 0.05312719085211144; "Q" = 3.982862030882444; "R" = 35.845758277942; "S" = 0.002075837785870323; "T" = 0.003934723291084594 }
    5 = @{ "E" = 3; "G" = 6.201571333333333; "H" = 18.604714; "I" = 0.05221490529364391; "J" = 0.07406232529850043; "K" = 3; "M" = 0.8924576666666667; "N" = 2.677373; "O" = 0.05524500589055249; "P" = 0.07382627543180768; "Q" = 5.534639881813556; "R" = 49.811758936322; "S" = 0.002884612750521998; "T" = 0.005467745626607231 }
    6 = @{ "E" = 3; "G" = 6.201571333333333; "H" = 18.604714; "I" = 0.05221490529364391; "J" = 0.07406232529850043; "K" = 2; "M" = 12.1977655; "N" = 24.395531; "O" = 0.7550673293176682; "P" = 0.6726859466018379; "Q" = 75.64531285552232; "R" = 453.8718771331339; "S" = 0.03942576909064668; "T" = 0.04982068540095501 }
    7 = @{ "E" = 3; "G" = 7.461641333333333; "H" = 22.384924; "I" = 0.06282422221945559; "J" = 0.0891107233935555; "K" = 3; "M" = 1.253707; "N" = 3.761121; "O" = 0.07760709912293902; "P" = 0.1037097015912075; "Q" = 9.354711971089332; "R" = 84.19240773980398; "S" = 0.004875605641106838; "T" = 0.009241646531722271 }
    8 = @{ "E" = 3; "G" = 7.461641333333333; "H" = 22.384924; "I" = 0.06282422221945559; "J" = 0.0891107233935555; "K" = 3; "M" = 1.168375666666667; "N" = 3.505127; "O" = 0.07232491018701337; "P" = 0.09665088552303537; "Q" = 8.718000167260888; "R" = 78.46200150534798; "S" = 0.004543756229591096; "T" = 0.008612630325585402 }
    9 = @{ "E" = 3; "G" = 7.461641333333333; "H" = 22.384924; "I" = 0.06282422221945559; "J" = 0.0891107233935555; "K" = 3; "M" = 0.6422343333333334; "N" = 1.926703; "O" = 0.03975565548182683; "P" = 0.05312719085211144; "Q" = 4.792122247285778; "R" = 43.129100225572; "S" = 0.002497618134470407; "T" = 0.004734202408699134 }
    10 = @{ "E" = 3; "G" = 7.461641333333333; "H" = 22.384924; "I" = 0.06282422221945559; "J" = 0.0891107233935555; "K" = 3; "M" = 0.8924576666666667; "N" = 2.677373; "O" = 0.05524500589055249; "P" = 0.07382627543180768; "Q" = 6.659199013850222; "R" = 59.932791124652; "S" = 0.003470724526583203; "T" = 0.006578712809180256 }
    11 = @{ "E" = 3; "G" = 7.461641333333333; "H" = 22.384924; "I" = 0.06282422221945559; "J" = 0.0891107233935555; "K" = 2; "M" = 12.1977655; "N" = 24.395531; "O" = 0.7550673293176682; "P" = 0.6726859466018379; "Q" = 91.01535122910732; "R" = 546.0921073746439; "S" = 0.04743651768770404; "T" = 0.05994353131836842 }
    12 = @{ "E" = 2; "G" = 105.106922; "H" = 210.213844; "I" = 0.8849608724869005; "J" = 0.836826951307944; "K" = 3; "M" = 1.253707; "N" = 3.761121; "O" = 0.07760709912293902; "P" = 0.1037097015912075; "Q" = 131.773283859854; "R" = 790.6397031591239; "S" = 0.06867924615101348; "T" = 0.08678707340362678 }
    13 = @{ "E" = 2; "G" = 105.106922; "H" = 210.213844; "I" = 0.8849608724869005; "J" = 0.836826951307944; "K" = 3; "M" = 1.168375666666667; "N" = 3.505127; "O" = 0.07232491018701337; "P" = 0.09665088552303537; "Q" = 122.8043700630313; "R" = 736.8262203781879; "S" = 0.06400471562163607; "T" = 0.0808800658734548 }
    14 = @{ "E" = 2; "G" = 105.106922; "H" = 210.213844; "I" = 0.8849608724869005; "J" = 0.836826951307944; "K" = 3; "M" = 0.6422343333333334; "N" = 1.926703; "O" = 0.03975565548182683; "P" = 0.05312719085211144; "Q" = 67.50327397938867; "R" = 405.019643876332; "S" = 0.0351821995614861; "T" = 0.04445826515232771 }
    15 = @{ "E" = 2; "G" = 105.106922; "H" = 210.213844; "I" = 0.8849608724869005; "J" = 0.836826951307944; "K" = 3; "M" = 0.8924576666666667; "N" = 2.677373; "O" = 0.05524500589055249; "P" = 0.07382627543180768; "Q" = 93.80347835863533; "R" = 562.8208701518121; "S" = 0.04888966861344728; "T" = 0.06177981699602019 }
    16 = @{ "E" = 2; "G" = 105.106922; "H" = 210.213844; "I" = 0.8849608724869005; "J" = 0.836826951307944; "K" = 2; "M" = 12.1977655; "N" = 24.395531; "O" = 0.7550673293176682; "P" = 0.6726859466018379; "Q" = 1282.069586982791; "R" = 5128.278347931164; "S" = 0.6682050425393175; "T" = 0.5629217298825144 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

Write-Host "Applied Natmi Hou advice update"
